$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'70.859.44"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = "  +2.54%  "
$ws.Cells.Item(3, 4).Value = "'3.569.62"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = "  +1.92%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).Value = "'614.18"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +6.31%  "
$ws.Cells.Item(6, 4).Value = "'172.62"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +0.88%  "
$ws.Cells.Item(7, 4).Value = "'0.620"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  +2.47%  "
$ws.Cells.Item(8, 4).Value = "'3.563.96"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  +1.89%  "
$ws.Cells.Item(9, 5).Value = "  -0.04%  "
$ws.Cells.Item(10, 5).Value = "  +5.13%  "
$ws.Cells.Item(11, 4).Value = "'7.28"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +12.83%  "
$ws.Cells.Item(12, 5).Value = "  +1.28%  "
$ws.Cells.Item(13, 4).Value = "'46.82"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +0.74%  "
$ws.Cells.Item(14, 4).Value = "'0.0000277"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +1.79%  "
$ws.Cells.Item(15, 4).Value = "'4.144.32"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "  +1.92%  "
$ws.Cells.Item(16, 5).Value = "  -1.17%  "
$ws.Cells.Item(17, 4).Value = "'618.96"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -0.08%  "
$ws.Cells.Item(18, 4).Value = "'3.570.18"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +2.04%  "
$ws.Cells.Item(19, 4).Value = "'70.951.98"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +2.77%  "
$ws.Cells.Item(20, 5).Value = "  -1.95%  "
$ws.Cells.Item(21, 5).Value = "  +0.55%  "
$ws.Cells.Item(22, 5).Value = "  +0.19%  "
$ws.Cells.Item(23, 4).Value = "'9.49"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -14.48%  "
$ws.Cells.Item(24, 4).Value = "'15.84"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -0.10%  "
$ws.Cells.Item(25, 4).Value = "'97.07"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +0.00%  "
$ws.Cells.Item(26, 5).Value = "  +1.51%  "
$ws.Cells.Item(27, 5).Value = "  -0.03%  "
$ws.Cells.Item(28, 5).Value = "  -0.10%  "
$ws.Cells.Item(29, 4).Value = "'33.71"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +3.46%  "
$ws.Cells.Item(30, 4).Value = "'9.11"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -1.93%  "
$ws.Cells.Item(31, 5).Value = "  +0.50%  "
$ws.Cells.Item(32, 5).Value = "  -2.23%  "
$ws.Cells.Item(33, 5).Value = "  -0.51%  "
$ws.Cells.Item(34, 4).Value = "'6.96"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -0.08%  "
$ws.Cells.Item(35, 4).Value = "'572.74"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -9.62%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "'0.101"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -0.98%  "
$ws.Cells.Item(37, 2).Value = "dogwifhat"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(37, 4).Value = "'3.62"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +5.62%  "
$ws.Cells.Item(38, 5).Value = "  +1.62%  "
$ws.Cells.Item(39, 4).Value = "'57.62"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +1.97%  "
$ws.Cells.Item(40, 4).Value = "'0.0474"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +6.16%  "
$ws.Cells.Item(41, 5).Value = "  +0.19%  "
$ws.Cells.Item(42, 5).Value = "  +4.98%  "
$ws.Cells.Item(43, 4).Value = "'3.372.19"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +0.37%  "
$ws.Cells.Item(44, 4).Value = "'0.321"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -1.71%  "
$ws.Cells.Item(45, 4).Value = "'2.99"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +8.22%  "
$ws.Cells.Item(46, 4).Value = "'33.07"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +1.08%  "
$ws.Cells.Item(47, 4).Value = "'0.0₃0705"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +2.32%  "
$ws.Cells.Item(48, 5).Value = "  +2.97%  "
$ws.Cells.Item(49, 5).Value = "  +0.80%  "
$ws.Cells.Item(50, 4).Value = "'133.77"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +1.24%  "
$ws.Cells.Item(51, 5).Value = "  +1.82%  "
